# Update the "Period" value cell (B6) to use the new dateTool.format(...) based
# expression instead of the old from.toString(...)/to.toString(...) expression.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'

# Update the position fix time cell (B9) to use the new dateTool.format(...) based
# expression instead of the old new("org.joda.time.DateTime", ...) expression.
$ws.Range("B9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", position.fixTime, locale, timezone)}'

# Move the active selection from G9 to B2, as recorded in the saved sheet view.
$ws.Range("B2").Select()
